# This script regenerates the "Password" column (E) for every existing user
# by re-running the create_user password routine, which joins First Name,
# Last Name and Class with randomly chosen separator characters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly generated passwords, keyed by the users row number on the "Users" sheet
$newPasswords = @{
    2 = 'Maximilian,De_Junious#2cn='
    3 = 'Ramon!Nunez_Gomez%3bn!'
    4 = 'Uelkue_Oemer=Uellaegoess=1an#'
    5 = 'Isis#Lanpher)4cn('
    6 = 'Maximilian%Galvin=2an='
    7 = 'Jannette,Laspina-4cn('
    8 = 'Marg)Dodich-3cn('
    9 = 'Carisa!Bannowsky!2an)'
    10 = 'David,Waisath&4cn='
    11 = 'Paulette^von_Reddig-Piette=2cn#'
    12 = 'Kirby!Latona.2an,'
    13 = 'Reed)Homewood(2bn!'
    14 = 'Blair^Pallafor-Zedian_5cn_'
    15 = 'Lon=Senemounnarat-Quillian#2cn-'
    16 = 'Vada%Isaac(2bn&'
    17 = 'Jeanett_Plancarte-4bn,'
    18 = 'Alex)Berteotti-Stirn(None&'
    19 = 'Robyn%Strycker,3an!'
    20 = 'Camille#Von_Verrill(3bn^'
    21 = 'Franz_Michael_Leopold&Deschner!4cn)'
    22 = 'Veola%Franzi%None.'
    23 = 'Chantelle=Cringle%3cn%'
    24 = 'Britney%Kosh!5bn='
    25 = 'Clayton^Derouchie%4bn%'
    26 = 'Beverlee&Doutt,5bn%'
    27 = 'Alma)Munley&1bn#'
    28 = 'Thad(Dornbos)5an,'
    29 = 'Arvilla&Mahala,2cn-'
    30 = 'Mirza!Ellingwood,None&'
    31 = 'Francie(de_Cardinalli-Sciola,3an&'
    32 = 'IRENEE_Gundry&None&'
    33 = 'IRENEE!Pinedo)3cn)'
    34 = 'Mirza,Pinsky#1bn#'
    35 = 'Francie-Pinsky1^1bn_'
    36 = 'Goldie=Pinsky2-1bn,'
    37 = 'A-nother#Pinsky3(1bn^'
}

foreach ($row in $newPasswords.Keys) {
    $ws.Cells.Item([int]$row, 5).Value = $newPasswords[$row]
}

Write-Output "Updated passwords for $($newPasswords.Count) users"
